$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ast_config")

# Update ast_condition (column M) status values to COMPLETE for rows 2-5
# Leading apostrophe preserves the existing "quote prefix" text cell style
$ws.Range("M2").Value = "'COMPLETE"
$ws.Range("M3").Value = "'COMPLETE"
$ws.Range("M4").Value = "'COMPLETE"
$ws.Range("M5").Value = "'COMPLETE"
